$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 814
$ws1.Range("F4").Value = 13576
$ws1.Range("F5").Value = 13385
$ws1.Range("F6").Value = 1037
$ws1.Range("F11").Value = 14
$ws1.Range("F12").Value = 25
$ws1.Range("F13").Value = 722
$ws1.Range("F15").Value = 49
$ws1.Range("F17").Value = 57
$ws1.Range("F20").Value = 417
$ws1.Range("F21").Value = 332
$ws1.Range("F23").Value = 477
$ws1.Range("F24").Value = 801
$ws1.Range("F25").Value = 57

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 42
$ws2.Range("F7").Value = 152
$ws2.Range("F8").Value = 937
$ws2.Range("F11").Value = 46

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 208
$ws3.Range("F3").Value = 83

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 208
$ws4.Range("F5").Value = 814
$ws4.Range("F6").Value = 13576
$ws4.Range("F7").Value = 13385
$ws4.Range("F8").Value = 1037
$ws4.Range("F13").Value = 14
$ws4.Range("F14").Value = 25
$ws4.Range("F15").Value = 722
$ws4.Range("F16").Value = 42
$ws4.Range("F19").Value = 49
$ws4.Range("F21").Value = 57
$ws4.Range("F26").Value = 83
$ws4.Range("F27").Value = 417
$ws4.Range("F28").Value = 332
$ws4.Range("F30").Value = 477
$ws4.Range("F31").Value = 801
$ws4.Range("F32").Value = 152
$ws4.Range("F33").Value = 937
$ws4.Range("F36").Value = 57
$ws4.Range("F37").Value = 46
